$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 12,20
$data[0,0] = "ECs"
$data[0,1] = "Lama1"
$data[0,2] = "Itga1"
$data[0,3] = "ECs"
$data[0,4] = 1
$data[0,5] = 0.3333333333333333
$data[0,6] = 0.02283333333333333
$data[0,7] = 0.0685
$data[0,8] = 0.05477178157813095
$data[0,9] = 0.05477178157813096
$data[0,10] = 3
$data[0,11] = 1
$data[0,12] = 75.59011833333334
$data[0,13] = 226.770355
$data[0,14] = 0.6588374259037486
$data[0,15] = 0.6588374259037486
$data[0,16] = 1.725974368611111
$data[0,17] = 15.5337693175
$data[0,18] = 0.03608569958709815
$data[0,19] = 0.03608569958709816
$data[1,0] = "ECs"
$data[1,1] = "Lama1"
$data[1,2] = "Itga1"
$data[1,3] = "FAPs"
$data[1,4] = 1
$data[1,5] = 0.3333333333333333
$data[1,6] = 0.02283333333333333
$data[1,7] = 0.0685
$data[1,8] = 0.05477178157813095
$data[1,9] = 0.05477178157813096
$data[1,10] = 3
$data[1,11] = 1
$data[1,12] = 12.15310033333333
$data[1,13] = 36.459301
$data[1,14] = 0.1059254505338229
$data[1,15] = 0.1059254505338229
$data[1,16] = 0.2774957909444444
$data[1,17] = 2.4974621185
$data[1,18] = 0.005801725640203664
$data[1,19] = 0.005801725640203664
$data[2,0] = "ECs"
$data[2,1] = "Lama1"
$data[2,2] = "Itga1"
$data[2,3] = "M2"
$data[2,4] = 1
$data[2,5] = 0.3333333333333333
$data[2,6] = 0.02283333333333333
$data[2,7] = 0.0685
$data[2,8] = 0.05477178157813095
$data[2,9] = 0.05477178157813096
$data[2,10] = 3
$data[2,11] = 1
$data[2,12] = 0.1727356666666667
$data[2,13] = 0.5182070000000001
$data[2,14] = 0.001505550255743542
$data[2,15] = 0.001505550255743542
$data[2,16] = 0.003944131055555557
$data[2,17] = 0.03549717950000001
$data[2,18] = 0.0000824616697624845
$data[2,19] = 0.00008246166976248451
$data[3,0] = "ECs"
$data[3,1] = "Lama1"
$data[3,2] = "Itga1"
$data[3,3] = "sCs"
$data[3,4] = 1
$data[3,5] = 0.3333333333333333
$data[3,6] = 0.02283333333333333
$data[3,7] = 0.0685
$data[3,8] = 0.05477178157813095
$data[3,9] = 0.05477178157813096
$data[3,10] = 3
$data[3,11] = 1
$data[3,12] = 26.81662666666667
$data[3,13] = 80.44988000000001
$data[3,14] = 0.233731573306685
$data[3,15] = 0.233731573306685
$data[3,16] = 0.6123129755555556
$data[3,17] = 5.510816780000001
$data[3,18] = 0.01280189468106665
$data[3,19] = 0.01280189468106665
$data[4,0] = "FAPs"
$data[4,1] = "Lama1"
$data[4,2] = "Itga1"
$data[4,3] = "ECs"
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 0.3395593333333333
$data[4,7] = 1.018678
$data[4,8] = 0.8145227578751427
$data[4,9] = 0.8145227578751428
$data[4,10] = 3
$data[4,11] = 1
$data[4,12] = 75.59011833333334
$data[4,13] = 226.770355
$data[4,14] = 0.6588374259037486
$data[4,15] = 0.6588374259037486
$data[4,16] = 25.66733018785445
$data[4,17] = 231.00597169069
$data[4,18] = 0.5366380771384813
$data[4,19] = 0.5366380771384814
$data[5,0] = "FAPs"
$data[5,1] = "Lama1"
$data[5,2] = "Itga1"
$data[5,3] = "FAPs"
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 0.3395593333333333
$data[5,7] = 1.018678
$data[5,8] = 0.8145227578751427
$data[5,9] = 0.8145227578751428
$data[5,10] = 3
$data[5,11] = 1
$data[5,12] = 12.15310033333333
$data[5,13] = 36.459301
$data[5,14] = 0.1059254505338229
$data[5,15] = 0.1059254505338229
$data[5,16] = 4.126698647119778
$data[5,17] = 37.14028782407799
$data[5,18] = 0.08627869009797647
$data[5,19] = 0.08627869009797647
$data[6,0] = "FAPs"
$data[6,1] = "Lama1"
$data[6,2] = "Itga1"
$data[6,3] = "M2"
$data[6,4] = 3
$data[6,5] = 1
$data[6,6] = 0.3395593333333333
$data[6,7] = 1.018678
$data[6,8] = 0.8145227578751427
$data[6,9] = 0.8145227578751428
$data[6,10] = 3
$data[6,11] = 1
$data[6,12] = 0.1727356666666667
$data[6,13] = 0.5182070000000001
$data[6,14] = 0.001505550255743542
$data[6,15] = 0.001505550255743542
$data[6,16] = 0.05865400781622224
$data[6,17] = 0.5278860703460001
$data[6,18] = 0.001226304946427857
$data[6,19] = 0.001226304946427857
$data[7,0] = "FAPs"
$data[7,1] = "Lama1"
$data[7,2] = "Itga1"
$data[7,3] = "sCs"
$data[7,4] = 3
$data[7,5] = 1
$data[7,6] = 0.3395593333333333
$data[7,7] = 1.018678
$data[7,8] = 0.8145227578751427
$data[7,9] = 0.8145227578751428
$data[7,10] = 3
$data[7,11] = 1
$data[7,12] = 26.81662666666667
$data[7,13] = 80.44988000000001
$data[7,14] = 0.233731573306685
$data[7,15] = 0.233731573306685
$data[7,16] = 9.105835873182222
$data[7,17] = 81.95252285864001
$data[7,18] = 0.1903796856922571
$data[7,19] = 0.1903796856922571
$data[8,0] = "sCs"
$data[8,1] = "Lama1"
$data[8,2] = "Itga1"
$data[8,3] = "ECs"
$data[8,4] = 3
$data[8,5] = 1
$data[8,6] = 0.05448866666666666
$data[8,7] = 0.163466
$data[8,8] = 0.1307054605467263
$data[8,9] = 0.1307054605467264
$data[8,10] = 3
$data[8,11] = 1
$data[8,12] = 75.59011833333334
$data[8,13] = 226.770355
$data[8,14] = 0.6588374259037486
$data[8,15] = 0.6588374259037486
$data[8,16] = 4.118804761158889
$data[8,17] = 37.06924285043
$data[8,18] = 0.08611364917816915
$data[8,19] = 0.08611364917816916
$data[9,0] = "sCs"
$data[9,1] = "Lama1"
$data[9,2] = "Itga1"
$data[9,3] = "FAPs"
$data[9,4] = 3
$data[9,5] = 1
$data[9,6] = 0.05448866666666666
$data[9,7] = 0.163466
$data[9,8] = 0.1307054605467263
$data[9,9] = 0.1307054605467264
$data[9,10] = 3
$data[9,11] = 1
$data[9,12] = 12.15310033333333
$data[9,13] = 36.459301
$data[9,14] = 0.1059254505338229
$data[9,15] = 0.1059254505338229
$data[9,16] = 0.6622062330295555
$data[9,17] = 5.959856097265999
$data[9,18] = 0.0138450347956428
$data[9,19] = 0.01384503479564281
$data[10,0] = "sCs"
$data[10,1] = "Lama1"
$data[10,2] = "Itga1"
$data[10,3] = "M2"
$data[10,4] = 3
$data[10,5] = 1
$data[10,6] = 0.05448866666666666
$data[10,7] = 0.163466
$data[10,8] = 0.1307054605467263
$data[10,9] = 0.1307054605467264
$data[10,10] = 3
$data[10,11] = 1
$data[10,12] = 0.1727356666666667
$data[10,13] = 0.5182070000000001
$data[10,14] = 0.001505550255743542
$data[10,15] = 0.001505550255743542
$data[10,16] = 0.009412136162444447
$data[10,17] = 0.08470922546200002
$data[10,18] = 0.0001967836395532013
$data[10,19] = 0.0001967836395532014
$data[11,0] = "sCs"
$data[11,1] = "Lama1"
$data[11,2] = "Itga1"
$data[11,3] = "sCs"
$data[11,4] = 3
$data[11,5] = 1
$data[11,6] = 0.05448866666666666
$data[11,7] = 0.163466
$data[11,8] = 0.1307054605467263
$data[11,9] = 0.1307054605467264
$data[11,10] = 3
$data[11,11] = 1
$data[11,12] = 26.81662666666667
$data[11,13] = 80.44988000000001
$data[11,14] = 0.233731573306685
$data[11,15] = 0.233731573306685
$data[11,16] = 1.461202231564444
$data[11,17] = 13.15082008408
$data[11,18] = 0.03054999293336118
$data[11,19] = 0.03054999293336119

$ws.Range("A2:T13").Value = $data
